$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Record_BuildingProduce")
$src.Copy([System.Reflection.Missing]::Value, $src)
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name $s.Index
}
